# Delete row 176 (the "カメラを止めろ" post) from Sheet1.
# This shifts all subsequent rows up by one, which matches the diff:
# dimension A1:C247 -> A1:C246, and all rows/cells below 176 renumbered down by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(176).Delete()
